$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "270.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.54%"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.32%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.722"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.06%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06132"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.16%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.748"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.49%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8555"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.47%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8982"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.28%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1429"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.89%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05070"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.14%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07156"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.95%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03189"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.62%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09044"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.14%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001528"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.13%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006086"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.21%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006026"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.90%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.12%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.183"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.28%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.278"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.59%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3088"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.57%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1282"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.15%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.839"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.56%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04241"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.31%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001179"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.08%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004148"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.58%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.02%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001679"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.83%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03971"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.43%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1119"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.62%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.42%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002040"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.62%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01178"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.80%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005140"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.68%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.03%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9071"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "433.55%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02995"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-16.62%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002003"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
